# Consolidate split text runs that share identical formatting back into a
# single run, mirroring the writer-side run-consolidation change described
# in the commit message.
#
# Slide 1 title: "Header" + " " + "with" + " " + [Consolas]"inline code"
#   -> merge the first four (identically-formatted, plain) runs into one
#      run "Header with ", leaving the differently-formatted "inline code"
#      run (Consolas) untouched.
# Slide 2 title: "Syntax" + " " + "highlighting" -> single run "Syntax highlighting"
# Slide 3 title: "Two" + " " + "column" + " " + "slide" -> single run "Two column slide"

$p = $ppt.ActivePresentation

# --- Slide 1: "Header with " stays merged, "inline code" run is untouched ---
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 12).Text = "Header with "

# --- Slide 2: whole title merges into one run ---
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, $tr2.Length).Text = "Syntax highlighting"

# --- Slide 3: whole title merges into one run ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, $tr3.Length).Text = "Two column slide"
